# "added the logo in svg"
#
# The deck currently has 2 slides:
#   1) A logo slide (picture + outline oval + "A" wordmark textbox)
#   2) Another logo slide (same picture/oval, different wordmark color/placement)
#
# The target adds a NEW logo slide (another wordmark variant) positioned
# right between the two existing slides. Its picture + oval are byte
# identical to slide 1's, so the most faithful way to create it through
# the PowerPoint object model is to duplicate slide 1 (which places the
# copy immediately after it, i.e. at position 2) and then restyle just
# the wordmark TextBox to match the new variant (position/size, word
# wrap, and font size).

$p = $ppt.ActivePresentation

$orig = $p.Slides.Item(1)
$newSlide = $orig.Duplicate().Item(1)

# The wordmark textbox is shape 3 (Picture 4, Oval 6, TextBox 8).
$textBox = $newSlide.Shapes.Item(3)

# Resize / reposition the textbox to the new variant's frame.
$textBox.Left = 335.9188232421875
$textBox.Top = 24.021024703979492
$textBox.Width = 238.1352081298828
$textBox.Height = 491.95782470703125

# New variant wraps text within the box instead of "none".
$textBox.TextFrame.WordWrap = $true

# New variant uses a smaller wordmark point size (400pt vs 450pt).
$textBox.TextFrame.TextRange.Font.Size = 400
